$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("students")

# Snapshot the current values of rows 2-5 (A:K) before overwriting anything,
# since the target state is a rearrangement of the existing rows' data.
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
$snapshot = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# New row order: row2<-old row4, row3<-old row5, row4<-old row2, row5<-old row3
$rowMap = @{ 2 = 4; 3 = 5; 4 = 2; 5 = 3 }

# The collaborator name on the row that becomes row 3 was inconsistently
# capitalized ("Manuel Villeda"); normalize it to "MANUEL VILLEDA" like the
# other occurrence in the sheet.
$snapshot[5]["J"] = "MANUEL VILLEDA"

# Columns that hold text in the source workbook (everything except the numeric
# "Cuenta Afectada" column G and the date "Fecha de Operacion" column I).
# These get written with a leading apostrophe (forcing text, which also keeps
# numeric-looking values like leading-zero IDs, and empty strings, as real
# text cells) and then have their style reset to Normal so no stray
# "quote prefix" number format gets attached to the cell.
$textColumns = @("A","B","C","D","E","F","H","J","K")

foreach ($newRow in 2..5) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $val = $src[$c]
        $cell = $ws.Range("$c$newRow")
        if ($textColumns -contains $c) {
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}

$wb.Save()
